# GreatOutdoor Delivery Time Report Sprint 2 Test Cases.xlsx
#
# The "Remarks" column (G) previously held a canned failure message
# ("Delivery Time Report - ":Invalid Argument Recieved) and the
# corresponding "Req: Reference" column (H) read "Fail" for test cases
# TC_01, TC_02, TC_03, TC_04, TC_06 and TC_09 (rows 5, 6, 7, 8, 10, 13).
# Those cases now pass, so the remark becomes a generic error note and
# the status flips to "Pass".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5, 6, 7, 8, 10, 13)
foreach ($r in $rows) {
    $ws.Range("G$r").Value = "Something Went wrong"
    $ws.Range("H$r").Value = "Pass"
}

# The workbook was last saved with cell G13 selected.
[void]$ws.Range("G13").Select()
